# [ADDITIONAL SCRAPING] added code to scrape more data about a player's
# batting performance in a match, also updated the excel sheets
#
# 1. Insert a new "Player Info" worksheet before the existing "ODI Batting"
#    sheet, with player metadata (ID, NAME, BATTING_HAND, BOWL_STYLE).
# 2. On the "ODI Batting" sheet, rename the MATCH_CARD_LINK column to
#    MATCH_CODE and replace the full scorecard URLs with just the bare
#    match code.

$wb = $excel.ActiveWorkbook

# --- 1. Create the new "Player Info" sheet -------------------------------
# Worksheets.Add() inserts the new sheet before the currently active sheet,
# which puts it first (ahead of "ODI Batting"), matching the target sheet
# order/sheetId/rId layout.
$playerInfo = $wb.Worksheets.Add()
$playerInfo.Name = "Player Info"

$playerInfo.Range("A1").Value = "ID"
$playerInfo.Range("B1").Value = "NAME"
$playerInfo.Range("C1").Value = "BATTING_HAND"
$playerInfo.Range("D1").Value = "BOWL_STYLE"

$idCell = $playerInfo.Range("A2")
$idCell.NumberFormat = "@"
$idCell.Value = "4992"
$playerInfo.Range("B2").Value = "Abid Ali"
$playerInfo.Range("C2").Value = "Right Handed"
$playerInfo.Range("D2").Value = "Right Arm Leg Break"

# Style the header row to match the workbook's existing header look
# (bold font, thin border, centered horizontally, top vertical alignment).
$headerRange = $playerInfo.Range("A1:D1")
$headerRange.Font.Bold = $true
$headerRange.HorizontalAlignment = -4108
$headerRange.VerticalAlignment = -4160
$headerRange.Borders.LineStyle = 1

# --- 2. Update the "ODI Batting" sheet ------------------------------------
$odi = $wb.Worksheets.Item("ODI Batting")

$odi.Range("D1").Value = "MATCH_CODE"

$matchCodes = @("4276", "4277", "4300", "4376", "4432", "4433")
for ($i = 0; $i -lt $matchCodes.Length; $i++) {
    $row = $i + 2
    $cell = $odi.Range("D" + $row)
    $cell.NumberFormat = "@"
    $cell.Value = $matchCodes[$i]
}
